# Insert a new team row ("Arran Katoko FC", id 19833277) into the
# worksheet, right above the "Pontaç0 F.C." row, shifting everything
# below it down by one row (and updating the hyperlinks accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 currently holds "Pontaç0 F.C." (id 20651178). Insert a fresh row
# above it so the existing data (rows 18-20) shifts down to rows 19-21.
$ws.Range("A18:C18").EntireRow.Insert()

# Fill in the new row with the new team's data.
$ws.Range("A18").Value = "Arran Katoko FC"
$ws.Range("B18").Value = 19833277
$ws.Range("C18").Value = "https://cartola.globo.com/#!/time/19833277"

# The row-insert does not shift the existing hyperlink anchors, so rebuild
# the whole hyperlink set against the final (post-insert) row numbers.
$ws.Hyperlinks.Delete()

$links = @(
    @{Row=2;  Id=32966},
    @{Row=3;  Id=184499},
    @{Row=4;  Id=186283},
    @{Row=5;  Id=287965},
    @{Row=6;  Id=1273719},
    @{Row=7;  Id=1326835},
    @{Row=8;  Id=1488983},
    @{Row=9;  Id=1747619},
    @{Row=10; Id=1867254},
    @{Row=11; Id=2371918},
    @{Row=12; Id=2916559},
    @{Row=13; Id=4088673},
    @{Row=14; Id=14709358},
    @{Row=15; Id=14933455},
    @{Row=16; Id=16411206},
    @{Row=17; Id=19209079},
    @{Row=18; Id=19833277},
    @{Row=19; Id=20651178},
    @{Row=20; Id=44810918},
    @{Row=21; Id=47775950}
)

foreach ($link in $links) {
    $cell = $ws.Cells.Item($link.Row, 3)
    $location = "!/time/" + $link.Id
    $ws.Hyperlinks.Add($cell, "https://cartola.globo.com/", $location) | Out-Null
    $cell.Style = "Hyperlink"
}
